# "kleine Änderungen an Folien" — split the "Nur mit manchen Boards..."
# sentence on the Arduino IDE 2.0 Debugging slide into three runs so the
# "SAMD MKR" wording becomes "SAMD21 basierende ... Boards wie MKR Zero)".

$p = $ppt.ActivePresentation

# Locate the slide that contains the placeholder text (robust to slide
# index, even though it is currently slide 8).
$targetSlide = $null
$targetShape = $null
foreach ($s in $p.Slides) {
    foreach ($shp in $s.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -like "*Nur mit manchen Boards unterst*") {
                $targetSlide = $s
                $targetShape = $shp
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Find the paragraph holding the sentence.
$para = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $candidate = $tr.Paragraphs($i)
    if ($candidate.Text -like "*Nur mit manchen Boards unterst*") {
        $para = $candidate
    }
}

# Replace the "SAMD MKR " fragment with "SAMD21 basierende " — this keeps
# the untouched prefix/suffix text intact while the middle portion becomes
# its own run, yielding:
#   "Nur mit manchen Boards unterstützt (" / "SAMD21 basierende " /
#   "Boards wie MKR Zero)"
$oldMid = "SAMD MKR "
$newMid = "SAMD21 basierende "

$fullText = $para.Text
$midStart = $fullText.IndexOf($oldMid) + 1   # 1-based for Characters()
$midLen = $oldMid.Length

$mid = $para.Characters($midStart, $midLen)
$mid.Text = $newMid
